$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1400
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1400
$ws.Range("N40").Value = -1750
$ws.Range("M40").ClearContents()

$ws.Range("H58").Value = 2296.2246
$ws.Range("I58").Value = 356.42856
$ws.Range("K58").Value = 1069.28568
$ws.Range("M58").Value = -919.28568

$ws.Range("H129").Value = 4033540.8
$ws.Range("J129").Value = 1321.4073
$ws.Range("L129").Value = 3964.2219
$ws.Range("N129").Value = -13964.2219

$ws.Range("H132").Value = 6062683
$ws.Range("I132").Value = 6898496
$ws.Range("J132").Value = 3037.5
$ws.Range("K132").Value = 20695488
$ws.Range("L132").Value = 9112.5
$ws.Range("M132").Value = -20692958
$ws.Range("N132").Value = -14172.5

$ws.Range("H135").Value = 1751.3529
$ws.Range("I135").Value = 1351.5333
$ws.Range("J135").Value = 4750
$ws.Range("K135").Value = 12163.7997
$ws.Range("L135").Value = 42750
$ws.Range("M135").Value = -9628.7997
$ws.Range("N135").Value = -47820

$ws.Range("H138").Value = 3457.8838
$ws.Range("I138").Value = 2217.3635
$ws.Range("J138").Value = 3884.3125
$ws.Range("K138").Value = 6652.0905
$ws.Range("L138").Value = 11652.9375
$ws.Range("M138").Value = -1512.0905
$ws.Range("N138").Value = -21932.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 15059
$ws.Range("J59").Value = 15059
$ws.Range("L59").Value = 15059
$ws.Range("N59").Value = -16667

$ws.Range("H74").Value = 970
$ws.Range("I74").Value = 902
$ws.Range("J74").Value = 1021
$ws.Range("K74").Value = 902
$ws.Range("L74").Value = 1021
$ws.Range("M74").Value = -28
$ws.Range("N74").Value = -2769

$ws.Range("H77").Value = 970
$ws.Range("I77").Value = 902
$ws.Range("J77").Value = 1021
$ws.Range("K77").Value = 4510
$ws.Range("L77").Value = 5105
$ws.Range("M77").Value = -142
$ws.Range("N77").Value = -13841

$ws.Range("H122").Value = 2560.842
$ws.Range("J122").Value = 3570.5
$ws.Range("L122").Value = 10711.5
$ws.Range("N122").Value = -15611.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2838.4614
$ws.Range("I99").Value = 1711.1111
$ws.Range("K99").Value = 1711.1111
$ws.Range("M99").Value = -213.1111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3202.0508
$ws.Range("I31").Value = 2106.9783
$ws.Range("J31").Value = 7076.923
$ws.Range("K31").Value = 2106.9783
$ws.Range("L31").Value = 7076.923
$ws.Range("M31").Value = -1811.9783
$ws.Range("N31").Value = -7666.923

$ws.Range("H34").Value = 3202.0508
$ws.Range("I34").Value = 2106.9783
$ws.Range("J34").Value = 7076.923
$ws.Range("K34").Value = 2106.9783
$ws.Range("L34").Value = 7076.923
$ws.Range("M34").Value = -1904.9783
$ws.Range("N34").Value = -7480.923

$ws.Range("H122").Value = 3089.6191
$ws.Range("I122").Value = 2636.7856
$ws.Range("K122").Value = 7910.3568
$ws.Range("M122").Value = -5460.3568

$ws.Range("H132").Value = 2862.6333
$ws.Range("I132").Value = 2173.95
$ws.Range("J132").Value = 4240
$ws.Range("K132").Value = 6521.849999999999
$ws.Range("L132").Value = 12720
$ws.Range("M132").Value = -3991.849999999999
$ws.Range("N132").Value = -17780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11409.866
$ws.Range("J87").Value = 15677.777
$ws.Range("L87").Value = 47033.331
$ws.Range("N87").Value = -49529.331

$ws.Range("H90").Value = 11409.866
$ws.Range("J90").Value = 15677.777
$ws.Range("L90").Value = 141099.993
$ws.Range("N90").Value = -153579.993

$ws.Range("H93").Value = 2450.923
$ws.Range("J93").Value = 2451.182
$ws.Range("L93").Value = 7353.545999999999
$ws.Range("N93").Value = -11097.546

$ws.Range("H120").Value = 18507.428
$ws.Range("I120").Value = 17388
$ws.Range("K120").Value = 52164
$ws.Range("M120").Value = -47326

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 61785.59
$ws.Range("I102").Value = 2571.4167
$ws.Range("K102").Value = 2571.4167
$ws.Range("M102").Value = -949.4167000000002

$ws.Range("H132").Value = 3306.238
$ws.Range("I132").Value = 3152.348
$ws.Range("K132").Value = 9457.044
$ws.Range("M132").Value = -6927.044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 83335210
$ws.Range("I22").Value = 250000780
$ws.Range("J22").Value = 2420.25
$ws.Range("K22").Value = 250000780
$ws.Range("L22").Value = 2420.25
$ws.Range("M22").Value = -250000485
$ws.Range("N22").Value = -3010.25

$ws.Range("H27").Value = 83335210
$ws.Range("I27").Value = 250000780
$ws.Range("J27").Value = 2420.25
$ws.Range("K27").Value = 250000780
$ws.Range("L27").Value = 2420.25
$ws.Range("M27").Value = -250000673
$ws.Range("N27").Value = -2634.25

$ws.Range("H51").Value = 24942
$ws.Range("J51").Value = 24942
$ws.Range("L51").Value = 24942
$ws.Range("N51").Value = -25898

$ws.Range("H122").Value = 3099.7932
$ws.Range("I122").Value = 2473.3684
$ws.Range("J122").Value = 4290
$ws.Range("K122").Value = 7420.1052
$ws.Range("L122").Value = 12870
$ws.Range("M122").Value = -4970.1052
$ws.Range("N122").Value = -17770

$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

$ws.Range("H132").Value = 4293.579
$ws.Range("I132").Value = 2840
$ws.Range("J132").Value = 4812.7144
$ws.Range("K132").Value = 8520
$ws.Range("L132").Value = 14438.1432
$ws.Range("M132").Value = -5990
$ws.Range("N132").Value = -19498.1432

$ws.Range("H136").Value = 3318.68
$ws.Range("I136").Value = 2966.6843
$ws.Range("J136").Value = 4433.3335
$ws.Range("K136").Value = 8900.052899999999
$ws.Range("L136").Value = 13300.0005
$ws.Range("M136").Value = -6350.052899999999
$ws.Range("N136").Value = -18400.0005

$ws.Range("H141").Value = 29833.334
$ws.Range("J141").Value = 29833.334
$ws.Range("L141").Value = 29833.334
$ws.Range("N141").Value = -40193.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1328.1154
$ws.Range("I107").Value = 573.7778
$ws.Range("K107").Value = 1721.3334
$ws.Range("M107").Value = 198.6666

$ws.Range("H122").Value = 1253874.2
$ws.Range("I122").Value = 2002400.8
$ws.Range("J122").Value = 6330
$ws.Range("K122").Value = 6007202.4
$ws.Range("L122").Value = 18990
$ws.Range("M122").Value = -6004752.4
$ws.Range("N122").Value = -23890

$ws.Range("H132").Value = 4494.364
$ws.Range("I132").Value = 1755.3182
$ws.Range("J132").Value = 7233.409
$ws.Range("K132").Value = 5265.9546
$ws.Range("L132").Value = 21700.227
$ws.Range("M132").Value = -2735.9546
$ws.Range("N132").Value = -26760.227
